# Update crypto price/volume figures per the Wed Jan 18 05:43:34 UTC 2023 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'303.49"
$ws.Range("D2").Style = $ws.Range("C2").Style
$ws.Range("E2").Value = "'1.44%"
$ws.Range("E2").Style = $ws.Range("C2").Style
$ws.Range("D3").Value = "'32.88"
$ws.Range("D3").Style = $ws.Range("C3").Style
$ws.Range("E3").Value = "'5.03%"
$ws.Range("E3").Style = $ws.Range("C3").Style
$ws.Range("D4").Value = "'4.947"
$ws.Range("D4").Style = $ws.Range("C4").Style
$ws.Range("E4").Value = "'-2.88%"
$ws.Range("E4").Style = $ws.Range("C4").Style
$ws.Range("D5").Value = "'0.07845"
$ws.Range("D5").Style = $ws.Range("C5").Style
$ws.Range("E5").Value = "'-1.48%"
$ws.Range("E5").Style = $ws.Range("C5").Style
$ws.Range("D6").Value = "'2.015"
$ws.Range("D6").Style = $ws.Range("C6").Style
$ws.Range("E6").Value = "'-11.54%"
$ws.Range("E6").Style = $ws.Range("C6").Style
$ws.Range("D7").Value = "'7.836"
$ws.Range("D7").Style = $ws.Range("C7").Style
$ws.Range("E7").Value = "'0.76%"
$ws.Range("E7").Style = $ws.Range("C7").Style
$ws.Range("D8").Value = "'3.807"
$ws.Range("D8").Style = $ws.Range("C8").Style
$ws.Range("E8").Value = "'-1.48%"
$ws.Range("E8").Style = $ws.Range("C8").Style
$ws.Range("E9").Value = "'-0.11%"
$ws.Range("E9").Style = $ws.Range("C9").Style
$ws.Range("D10").Value = "'0.1756"
$ws.Range("D10").Style = $ws.Range("C10").Style
$ws.Range("E10").Value = "'0.98%"
$ws.Range("E10").Style = $ws.Range("C10").Style
$ws.Range("D11").Value = "'0.07851"
$ws.Range("D11").Style = $ws.Range("C11").Style
$ws.Range("E11").Value = "'4.15%"
$ws.Range("E11").Style = $ws.Range("C11").Style
$ws.Range("D12").Value = "'0.08684"
$ws.Range("D12").Style = $ws.Range("C12").Style
$ws.Range("E12").Value = "'-6.29%"
$ws.Range("E12").Style = $ws.Range("C12").Style
$ws.Range("D13").Value = "'0.03147"
$ws.Range("D13").Style = $ws.Range("C13").Style
$ws.Range("E13").Value = "'3.51%"
$ws.Range("E13").Style = $ws.Range("C13").Style
$ws.Range("D14").Value = "'0.1005"
$ws.Range("D14").Style = $ws.Range("C14").Style
$ws.Range("E14").Value = "'0.08%"
$ws.Range("E14").Style = $ws.Range("C14").Style
$ws.Range("D15").Value = "'0.001517"
$ws.Range("D15").Style = $ws.Range("C15").Style
$ws.Range("E15").Value = "'0.42%"
$ws.Range("E15").Style = $ws.Range("C15").Style
$ws.Range("D16").Value = "'0.005876"
$ws.Range("D16").Style = $ws.Range("C16").Style
$ws.Range("E16").Value = "'-0.20%"
$ws.Range("E16").Style = $ws.Range("C16").Style
$ws.Range("E17").Value = "'-0.58%"
$ws.Range("E17").Style = $ws.Range("C17").Style
$ws.Range("D18").Value = "'2.154"
$ws.Range("D18").Style = $ws.Range("C18").Style
$ws.Range("D19").Value = "'0.3308"
$ws.Range("D19").Style = $ws.Range("C19").Style
$ws.Range("E19").Value = "'1.11%"
$ws.Range("E19").Style = $ws.Range("C19").Style
$ws.Range("D20").Value = "'0.1318"
$ws.Range("D20").Style = $ws.Range("C20").Style
$ws.Range("E20").Value = "'-1.18%"
$ws.Range("E20").Style = $ws.Range("C20").Style
$ws.Range("D21").Value = "'4.314"
$ws.Range("D21").Style = $ws.Range("C21").Style
$ws.Range("E21").Value = "'10.18%"
$ws.Range("E21").Style = $ws.Range("C21").Style
$ws.Range("D22").Value = "'0.1991"
$ws.Range("D22").Style = $ws.Range("C22").Style
$ws.Range("E22").Value = "'17.05%"
$ws.Range("E22").Style = $ws.Range("C22").Style
$ws.Range("D23").Value = "'0.04562"
$ws.Range("D23").Style = $ws.Range("C23").Style
$ws.Range("E23").Value = "'-1.04%"
$ws.Range("E23").Style = $ws.Range("C23").Style
$ws.Range("D24").Value = "'0.001224"
$ws.Range("D24").Style = $ws.Range("C24").Style
$ws.Range("E24").Value = "'-2.07%"
$ws.Range("E24").Style = $ws.Range("C24").Style
$ws.Range("E25").Value = "'-0.81%"
$ws.Range("E25").Style = $ws.Range("C25").Style
$ws.Range("E26").Value = "'4.20%"
$ws.Range("E26").Style = $ws.Range("C26").Style
$ws.Range("D39").Value = "'0.01737"
$ws.Range("D39").Style = $ws.Range("C39").Style
$ws.Range("E39").Value = "'-1.30%"
$ws.Range("E39").Style = $ws.Range("C39").Style
$ws.Range("D40").Value = "'0.04791"
$ws.Range("D40").Style = $ws.Range("C40").Style
$ws.Range("E40").Value = "'3.36%"
$ws.Range("E40").Style = $ws.Range("C40").Style
$ws.Range("D41").Value = "'0.007496"
$ws.Range("D41").Style = $ws.Range("C41").Style
$ws.Range("E41").Value = "'7.60%"
$ws.Range("E41").Style = $ws.Range("C41").Style
$ws.Range("E42").Value = "'-0.21%"
$ws.Range("E42").Style = $ws.Range("C42").Style
$ws.Range("D43").Value = "'0.002340"
$ws.Range("D43").Style = $ws.Range("C43").Style
$ws.Range("E43").Value = "'6.87%"
$ws.Range("E43").Style = $ws.Range("C43").Style
$ws.Range("E44").Value = "'2.75%"
$ws.Range("E44").Style = $ws.Range("C44").Style
$ws.Range("D45").Value = "'0.00006259"
$ws.Range("D45").Style = $ws.Range("C45").Style
$ws.Range("E45").Value = "'-0.56%"
$ws.Range("E45").Style = $ws.Range("C45").Style
$ws.Range("E46").Value = "'0.09%"
$ws.Range("E46").Style = $ws.Range("C46").Style
$ws.Range("E47").Value = "'-61.10%"
$ws.Range("E47").Style = $ws.Range("C47").Style
$ws.Range("D48").Value = "'0.8234"
$ws.Range("D48").Style = $ws.Range("C48").Style
$ws.Range("E48").Value = "'10.26%"
$ws.Range("E48").Style = $ws.Range("C48").Style
$ws.Range("E49").Value = "'0.09%"
$ws.Range("E49").Style = $ws.Range("C49").Style
$ws.Range("E50").Value = "'0.09%"
$ws.Range("E50").Style = $ws.Range("C50").Style
